$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the whole "Pinterest" source paragraph:
#    "https://no.pinterest.com/pin/647040671446427690/ Bakgrunnsbilde til
#    headlinen." — find it by its hyperlink text and delete the entire
#    paragraph (including its paragraph mark) so the following sources shift
#    up by one paragraph, exactly like the diff shows.
# ---------------------------------------------------------------------------
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*no.pinterest.com*") {
        $targetPara = $p
        break
    }
}
if ($targetPara -ne $null) {
    $targetPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) Append two new "Image by ..." Pixabay credit paragraphs after the
#    "... Bakgrunsbilde til about us." paragraph (right before the document's
#    final empty paragraph), preceded by one new blank paragraph.
#    Built via InsertXML on the trailing empty paragraph so we can reproduce
#    the exact run / proofErr structure from the target OOXML.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newXml = '<w:body ' + $w + '>' `
  + '<w:p/>' `
  + '<w:p>' `
    + '<w:r><w:t xml:space="preserve">Image by &lt;a href="https://pixabay.com/users/bluebudgie-4333174/?utm_source=link-attribution&amp;utm_medium=referral&amp;utm_campaign=image&amp;utm_content=2735173"&gt;Davie </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>Bicker</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t>&lt;/a&gt; from &lt;a href="https://pixabay.com//?utm_source=link-attribution&amp;utm_medium=referral&amp;utm_campaign=image&amp;utm_content=2735173"&gt;Pixabay&lt;/a&gt;</w:t></w:r>' `
  + '</w:p>' `
  + '<w:p>' `
    + '<w:r><w:t xml:space="preserve">Image by &lt;a href="https://pixabay.com/users/lockenkopf-4101190/?utm_source=link-attribution&amp;utm_medium=referral&amp;utm_campaign=image&amp;utm_content=2817112"&gt;Tanja </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>Schulte</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t>&lt;/a&gt; from &lt;a href="https://pixabay.com//?utm_source=link-attribution&amp;utm_medium=referral&amp;utm_campaign=image&amp;utm_content=2817112"&gt;Pixabay&lt;/a&gt;</w:t></w:r>' `
  + '</w:p>' `
  + '<w:p/>' `
  + '</w:body>'

$lastPara.Range.InsertXML($newXml)
